# Apply the weekly shuffle of Fecha/Calidad/Volumen/Precios/Origen across rows 2-27
# (row 23 is unchanged - it maps to itself in the permutation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 25 (pre-edit snapshot)
$ws.Range("D2").Value2 = 44232
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 16000
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 3000
$ws.Range("O2").Value = "Provincia de Chacabuco"
$ws.Range("P2").Value = 30

# Row 3 <- source row 18 (pre-edit snapshot)
$ws.Range("D3").Value2 = 44167
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 7000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("O3").Value = "Provincia de Chacabuco"
$ws.Range("P3").Value = 30

# Row 4 <- source row 9 (pre-edit snapshot)
$ws.Range("D4").Value2 = 44204
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 7000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = "Provincia de Chacabuco"
$ws.Range("P4").Value = 30

# Row 5 <- source row 24 (pre-edit snapshot)
$ws.Range("D5").Value2 = 44168
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 7000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("O5").Value = "Provincia de Chacabuco"
$ws.Range("P5").Value = 30

# Row 6 <- source row 27 (pre-edit snapshot)
$ws.Range("D6").Value2 = 44189
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 16000
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 3000
$ws.Range("O6").Value = "Provincia de Chacabuco"
$ws.Range("P6").Value = 30

# Row 7 <- source row 19 (pre-edit snapshot)
$ws.Range("D7").Value2 = 44230
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 16000
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("O7").Value = "Provincia de Chacabuco"
$ws.Range("P7").Value = 30

# Row 8 <- source row 17 (pre-edit snapshot)
$ws.Range("D8").Value2 = 44229
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 16000
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("O8").Value = "Provincia de Chacabuco"
$ws.Range("P8").Value = 30

# Row 9 <- source row 12 (pre-edit snapshot)
$ws.Range("D9").Value2 = 44210
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 8800
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 2750
$ws.Range("O9").Value = "Provincia de Chacabuco"
$ws.Range("P9").Value = 28

# Row 10 <- source row 15 (pre-edit snapshot)
$ws.Range("D10").Value2 = 44166
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 7000
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("O10").Value = "Provincia de Chacabuco"
$ws.Range("P10").Value = 30

# Row 11 <- source row 10 (pre-edit snapshot)
$ws.Range("D11").Value2 = 44209
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 7000
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 2750
$ws.Range("O11").Value = "Provincia de Chacabuco"
$ws.Range("P11").Value = 28

# Row 12 <- source row 2 (pre-edit snapshot)
$ws.Range("D12").Value2 = 44161
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 7000
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("O12").Value = "Provincia de Chacabuco"
$ws.Range("P12").Value = 30

# Row 13 <- source row 16 (pre-edit snapshot)
$ws.Range("D13").Value2 = 44181
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 12000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("O13").Value = "Provincia de Chacabuco"
$ws.Range("P13").Value = 30

# Row 14 <- source row 6 (pre-edit snapshot)
$ws.Range("D14").Value2 = 44186
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = 3000
$ws.Range("O14").Value = "Provincia de Chacabuco"
$ws.Range("P14").Value = 30

# Row 15 <- source row 20 (pre-edit snapshot)
$ws.Range("D15").Value2 = 44245
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 9000
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 3000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 30

# Row 16 <- source row 21 (pre-edit snapshot)
$ws.Range("D16").Value2 = 44245
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 25

# Row 17 <- source row 4 (pre-edit snapshot)
$ws.Range("D17").Value2 = 44188
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 12000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 3000
$ws.Range("O17").Value = "Provincia de Chacabuco"
$ws.Range("P17").Value = 30

# Row 18 <- source row 8 (pre-edit snapshot)
$ws.Range("D18").Value2 = 44160
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("O18").Value = "Provincia de Chacabuco"
$ws.Range("P18").Value = 30

# Row 19 <- source row 3 (pre-edit snapshot)
$ws.Range("D19").Value2 = 44214
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 7000
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("O19").Value = "Provincia de Chacabuco"
$ws.Range("P19").Value = 30

# Row 20 <- source row 11 (pre-edit snapshot)
$ws.Range("D20").Value2 = 44159
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 7000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 3000
$ws.Range("O20").Value = "Provincia de Chacabuco"
$ws.Range("P20").Value = 30

# Row 21 <- source row 7 (pre-edit snapshot)
$ws.Range("D21").Value2 = 44187
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 12000
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 3000
$ws.Range("O21").Value = "Provincia de Chacabuco"
$ws.Range("P21").Value = 30

# Row 22 <- source row 26 (pre-edit snapshot)
$ws.Range("D22").Value2 = 44215
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 16000
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("O22").Value = "Provincia de Chacabuco"
$ws.Range("P22").Value = 30

# Row 24 <- source row 5 (pre-edit snapshot)
$ws.Range("D24").Value2 = 44600
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 1300
$ws.Range("K24").Value = 3500
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 3808
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 38

# Row 25 <- source row 22 (pre-edit snapshot)
$ws.Range("D25").Value2 = 44231
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 12000
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = 3000
$ws.Range("O25").Value = "Provincia de Chacabuco"
$ws.Range("P25").Value = 30

# Row 26 <- source row 13 (pre-edit snapshot)
$ws.Range("D26").Value2 = 44602
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 12000
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = 3000
$ws.Range("O26").Value = "Provincia de Chacabuco"
$ws.Range("P26").Value = 30

# Row 27 <- source row 14 (pre-edit snapshot)
$ws.Range("D27").Value2 = 44602
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = 2500
$ws.Range("O27").Value = "Provincia de Chacabuco"
$ws.Range("P27").Value = 25

